# Fruta / hortaliza, semanal
# Insert 3 new rows at the top of the current week's block (rows 742-744),
# pushing the existing rows 742..832 down to 745..835, and populate the
# 3 new rows with this week's freshly scraped price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by 3 rows (Excel default Insert = shift down,
# copying format from the row above into the freshly inserted rows).
$ws.Rows("742:744").Insert()

function Set-Row($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value  = 11
    $ws.Cells.Item($Row, 2).Value  = "Vega Monumental Concepción"
    $ws.Cells.Item($Row, 3).Value  = "Bíobío"
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = 8
    $ws.Cells.Item($Row, 6).Value  = "Fruta"
    $ws.Cells.Item($Row, 7).Value  = 100108
    $ws.Cells.Item($Row, 8).Value  = "Tropicales y subtropicales"
    $ws.Cells.Item($Row, 9).Value  = 100108006
    $ws.Cells.Item($Row, 10).Value = "Plátano"
    $ws.Cells.Item($Row, 11).Value = "Sin especificar"
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = "$/caja 20 kilos"
    $ws.Cells.Item($Row, 18).Value = "Ecuador"
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = 20
}

Set-Row 742 45154 "Maduro"         100 11000 11000 11000 550
Set-Row 743 45154 "Pintón"         300 12000 12000 12000 600
Set-Row 744 45154 "Primera Pintón" 300 15000 15000 15000 750
